# Apply latest crypto price/volume scrape results (GitHub Actions data refresh).
# For each changed row, update only the cells that actually moved: Coin/Link
# text in B/C (rows 47-48 swapped ranking), Price in D, Volume(1h) in E.
# Numeric-looking Price strings are written with a leading apostrophe so Excel
# keeps them as text (matching the original inlineStr cells) instead of coercing
# them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.177.18'
$ws.Range('E2').Value = '  -0.44%  '
# Row 3
$ws.Range('D3').Value = '1.853.62'
# Row 4
$ws.Range('E4').Value = '  -0.04%  '
# Row 5
$ws.Range('D5').Value = '''235.31'
$ws.Range('E5').Value = '  -0.11%  '
# Row 6
$ws.Range('D6').Value = '''0.9999'
$ws.Range('E6').Value = '  -0.03%  '
# Row 7
$ws.Range('D7').Value = '''0.4691'
$ws.Range('E7').Value = '  +0.17%  '
# Row 8
$ws.Range('D8').Value = '''0.2888'
$ws.Range('E8').Value = '  +1.62%  '
# Row 9
$ws.Range('D9').Value = '''0.06553'
$ws.Range('E9').Value = '  +0.32%  '
# Row 10
$ws.Range('D10').Value = '''21.81'
$ws.Range('E10').Value = '  +2.27%  '
# Row 11
$ws.Range('D11').Value = '''0.07971'
$ws.Range('E11').Value = '  +1.32%  '
# Row 12
$ws.Range('D12').Value = '''97.44'
$ws.Range('E12').Value = '  -0.49%  '
# Row 13
$ws.Range('D13').Value = '1.854.55'
$ws.Range('E13').Value = '  -0.75%  '
# Row 14
$ws.Range('D14').Value = '''5.105'
$ws.Range('E14').Value = '  +0.09%  '
# Row 15
$ws.Range('D15').Value = '''0.6758'
$ws.Range('E15').Value = '  +0.15%  '
# Row 16
$ws.Range('D16').Value = '''268.71'
$ws.Range('E16').Value = '  -3.75%  '
# Row 17
$ws.Range('D17').Value = '30.150.76'
$ws.Range('E17').Value = '  -0.52%  '
# Row 18
$ws.Range('D18').Value = '''13.61'
$ws.Range('E18').Value = '  +6.99%  '
# Row 19
$ws.Range('D19').Value = '''0.000007706'
$ws.Range('E19').Value = '  +5.68%  '
# Row 20
$ws.Range('E20').Value = '  -0.03%  '
# Row 21
$ws.Range('D21').Value = '2.098.78'
$ws.Range('E21').Value = '  -0.77%  '
# Row 22
$ws.Range('E22').Value = '  -0.04%  '
# Row 23
$ws.Range('D23').Value = '''5.193'
$ws.Range('E23').Value = '  -5.69%  '
# Row 24
$ws.Range('D24').Value = '''6.137'
$ws.Range('E24').Value = '  -0.42%  '
# Row 25
$ws.Range('D25').Value = '''167.14'
# Row 26
$ws.Range('D26').Value = '''9.160'
$ws.Range('E26').Value = '  -0.07%  '
# Row 27
$ws.Range('E27').Value = '  -1.34%  '
# Row 28
$ws.Range('E28').Value = '  -0.15%  '
# Row 29
$ws.Range('E29').Value = '  +0.09%  '
# Row 30
$ws.Range('D30').Value = '''0.09856'
$ws.Range('E30').Value = '  +2.28%  '
# Row 31
$ws.Range('D31').Value = '''1.463'
$ws.Range('E31').Value = '  -0.99%  '
# Row 32
$ws.Range('D32').Value = '''4.285'
$ws.Range('E32').Value = '  -1.96%  '
# Row 33
$ws.Range('D33').Value = '''3.995'
$ws.Range('E33').Value = '  -2.45%  '
# Row 35
$ws.Range('E35').Value = '  -0.83%  '
# Row 36
$ws.Range('D36').Value = '''0.6984'
$ws.Range('E36').Value = '  -1.14%  '
# Row 37
$ws.Range('D37').Value = '''2.705'
$ws.Range('E37').Value = '  -0.61%  '
# Row 38
$ws.Range('D38').Value = '''0.01870'
$ws.Range('E38').Value = '  +0.71%  '
# Row 39
$ws.Range('D39').Value = '''2.602'
$ws.Range('E39').Value = '  +2.90%  '
# Row 40
$ws.Range('D40').Value = '''6.321'
$ws.Range('E40').Value = '  +0.69%  '
# Row 41
$ws.Range('E41').Value = '  -1.23%  '
# Row 42
$ws.Range('D42').Value = '''1.932'
$ws.Range('E42').Value = '  -0.57%  '
# Row 43
$ws.Range('D43').Value = '''0.9990'
$ws.Range('E43').Value = '  -0.12%  '
# Row 44
$ws.Range('D44').Value = '''0.8389'
$ws.Range('E44').Value = '  -1.30%  '
# Row 45
$ws.Range('D45').Value = '''103.10'
# Row 46
$ws.Range('D46').Value = '''0.4129'
$ws.Range('E46').Value = '  -1.25%  '
# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''9.174'
$ws.Range('E47').Value = '  -0.84%  '
# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '''936.35'
$ws.Range('E48').Value = '  +0.03%  '
# Row 49
$ws.Range('D49').Value = '''7.016'
$ws.Range('E49').Value = '  -2.10%  '
# Row 50
$ws.Range('D50').Value = '''33.83'
$ws.Range('E50').Value = '  -0.98%  '
# Row 51
$ws.Range('E51').Value = '  +0.34%  '
